$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from its current location (end of the 3rd
#    paragraph, just before "The following sections...") to the very start
#    of the document (immediately before the first run of paragraph 1).
# ---------------------------------------------------------------------------

# Remove the existing (collapsed) "_GoBack" bookmark, wherever it is.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Word's Bookmarks.Add() has trouble anchoring a brand-new, zero-length
# bookmark exactly at absolute character position 0 (it silently snaps to
# cover the whole first paragraph instead of staying collapsed). To avoid
# that - and to avoid disturbing the existing runs of paragraph 1 (which
# must remain untouched) - insert a temporary paragraph in front of
# everything, anchor the bookmark at the (non-zero) start of what is now
# the *second* paragraph (i.e. the original first paragraph), and then
# delete the temporary paragraph again. Deleting a whole separate paragraph
# does not cause Word to re-merge the runs that make up paragraph 1, so
# their original run boundaries/formatting survive intact, while the
# bookmark correctly slides back down to character position 0.
$d.Range(0, 0).InsertParagraphBefore()

$origPara1Start = $d.Paragraphs.Item(2).Range.Start
$anchor = $d.Range($origPara1Start, $origPara1Start)
$d.Bookmarks.Add("_GoBack", $anchor)

$d.Paragraphs.Item(1).Range.Delete()

# ---------------------------------------------------------------------------
# 2. Merge the three runs that spell out the S3 bucket name (with the
#    spell-check proofErr wrapper around "dmanwill") into a single run.
# ---------------------------------------------------------------------------

$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$bucketPhrase = $openQuote + "dmanwill-project-dataset" + $closeQuote

$d.Content.Find.Execute($bucketPhrase, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $bucketPhrase, 2) | Out-Null
